$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original (default) style for column D, then force the cells to
# be treated as text so numeric-looking price strings are not converted to
# floating point numbers. The style is restored afterwards so no visible
# formatting change is left behind.
$origStyle = $ws.Range("D2").Style
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range('D2').Value = '67.450.16'
$ws.Range('D3').Value = '3.509.15'
$ws.Range('D5').Value = '199.15'
$ws.Range('D6').Value = '551.15'
$ws.Range('D7').Value = '0.622'
$ws.Range('D8').Value = '3.496.65'
$ws.Range('D11').Value = '61.62'
$ws.Range('D13').Value = '0.0000268'
$ws.Range('D14').Value = '9.79'
$ws.Range('D15').Value = '4.076.66'
$ws.Range('D16').Value = '3.517.70'
$ws.Range('D18').Value = '67.189.03'
$ws.Range('D19').Value = '18.32'
$ws.Range('D20').Value = '11.80'
$ws.Range('D22').Value = '392.08'
$ws.Range('D23').Value = '3.97'
$ws.Range('D24').Value = '11.85'
$ws.Range('D25').Value = '84.66'
$ws.Range('D26').Value = '3.86'
$ws.Range('D27').Value = '12.21'
$ws.Range('D28').Value = '2.81'
$ws.Range('D29').Value = '8.81'
$ws.Range('D30').Value = '710.56'
$ws.Range('D31').Value = '31.00'
$ws.Range('D32').Value = '6.99'
$ws.Range('D33').Value = '11.68'
$ws.Range('D34').Value = '63.88'
$ws.Range('D35').Value = '0.110'
$ws.Range('D36').Value = '38.33'
$ws.Range('D38').Value = '0.393'
$ws.Range('D39').Value = '3.01'
$ws.Range('D40').Value = '0.999'
$ws.Range('D41').Value = '0.130'
$ws.Range('D42').Value = '3.063.41'
$ws.Range('D43').Value = '0.0₃0678'
$ws.Range('D44').Value = '2.79'
$ws.Range('D45').Value = '2.50'
$ws.Range('D46').Value = '0.0406'
$ws.Range('D47').Value = '0.129'
$ws.Range('D48').Value = '2.58'
$ws.Range('D49').Value = '138.05'
$ws.Range('D50').Value = '8.23'

# Restore the original cell style now that the text values are set.
$dRange.Style = $origStyle

# --- Coin name / link / volume (columns B, C, E) updates ---
$ws.Range('E2').Value = '  -1.57%  '
$ws.Range('E3').Value = '  -3.89%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('E6').Value = '  -4.69%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').Value = '  -4.15%  '
$ws.Range('E11').Value = '  +8.66%  '
$ws.Range('E12').Value = '  -8.03%  '
$ws.Range('E13').Value = '  -9.47%  '
$ws.Range('E14').Value = '  -3.48%  '
$ws.Range('E15').Value = '  -3.79%  '
$ws.Range('E16').Value = '  -3.52%  '
$ws.Range('E17').Value = '  -2.03%  '
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('E20').Value = '  -6.41%  '
$ws.Range('E21').Value = '  -6.63%  '
$ws.Range('E22').Value = '  -3.29%  '
$ws.Range('E23').Value = '  -7.13%  '
$ws.Range('E24').Value = '  -9.93%  '
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('E26').Value = '  -0.28%  '
$ws.Range('E27').Value = '  -3.88%  '
$ws.Range('E28').Value = '  -5.53%  '
$ws.Range('E29').Value = '  -4.93%  '
$ws.Range('E30').Value = '  +2.04%  '
$ws.Range('E31').Value = '  -3.44%  '
$ws.Range('E32').Value = '  -15.72%  '
$ws.Range('E33').Value = '  -5.13%  '
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('E35').Value = '  -6.11%  '
$ws.Range('E36').Value = '  -10.65%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -8.57%  '
$ws.Range('E39').Value = '  -4.90%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E41').Value = '  -8.74%  '
$ws.Range('E42').Value = '  -4.63%  '
$ws.Range('E43').Value = '  -14.88%  '
$ws.Range('E44').Value = '  +5.75%  '
$ws.Range('E45').Value = '  -14.38%  '
$ws.Range('E46').Value = '  -4.35%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E47').Value = '  -3.06%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E48').Value = '  -14.30%  '
$ws.Range('E50').Value = '  -8.24%  '
$ws.Range('E51').Value = '  -6.90%  '

Write-Output "Applied 99 cell updates"
